$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.157.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.974.61'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.07%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.971.41'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.64%  '
$ws.Range("E10").Value = '  -4.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.77'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.453'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.41%  '
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.463.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.01'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.094.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.969.38'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.682'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -8.13%  '
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.55%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.53%  '
$ws.Range("E32").Value = '  -5.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0793'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.26'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("E39").Value = '  -4.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.12'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.119'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.78'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '390.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.40%  '
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.264'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.21%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.693.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '37.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("E50").Value = '  -0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.32%  '
